# Fix contact information missing from short resumes:
# insert a centered contact-info paragraph right after the name/title
# paragraph ("Dheeraj Chand"), matching the long-resume layout.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Dheeraj Chand",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX",
    2
)
